$d = $word.ActiveDocument

# Locate the paragraph that starts with "Tracing the Evolution" (the title-idea
# paragraph split across 3 runs: "Tracing the Evolution" / " of Biological
# signals" / ": The case studies of Vision and Chemokines") and the paragraph
# right after the following blank line that reads "Tracing Evolutionary
# Pathways" + ": The case studies of Vision and Chemokines".
#
# We replace that whole span - the two content paragraphs plus the blank
# paragraph between them - with:
#   1) the original "Tracing the Evolution ..." paragraph, but now with its
#      tail merged into a single run ending in "...Vision and Chemokines"
#   2) a blank paragraph
#   3) a (new) copy of the original "Tracing Evolutionary Pathways: The case
#      studies of Vision and Chemokines" paragraph
#   4) a blank paragraph
#   5) a reworded "Tracing Evolutionary Pathways: Insights from Vision and
#      Chemokines" paragraph (split into 3 runs)

$startPara = $null
$midEndPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t.StartsWith("Tracing the Evolution of Biological signals")) {
        $startPara = $p
    }
    if ($midEndPara -eq $null -and $t.StartsWith("Tracing Evolutionary Pathways")) {
        $midEndPara = $p
    }
}
# The paragraph immediately after "Tracing Evolutionary Pathways: ..." is a
# blank paragraph; it gets absorbed into the replacement below (one of the
# three trailing blank paragraphs is consumed so the net paragraph count
# only grows by one).
$trailingBlank = $d.Paragraphs.Item($midEndPara.Index + 1)

$target = $d.Range($startPara.Range.Start, $trailingBlank.Range.End)

$xml = '<w:p><w:r><w:t>Tracing the Evolution</w:t></w:r><w:r><w:t xml:space="preserve"> of Biological signals: The case studies of Vision and Chemokines</w:t></w:r></w:p>' `
     + '<w:p/>' `
     + '<w:p><w:r><w:t>Tracing Evolutionary Pathways</w:t></w:r><w:r><w:t>: The case studies of Vision and Chemokines</w:t></w:r></w:p>' `
     + '<w:p/>' `
     + '<w:p><w:r><w:t xml:space="preserve">Tracing Evolutionary Pathways: </w:t></w:r><w:r><w:t>Insights from</w:t></w:r><w:r><w:t xml:space="preserve"> Vision and Chemokines</w:t></w:r></w:p>'

$target.InsertXML($xml) | Out-Null

Write-Output $d.Paragraphs.Count
